# chore: update Sheets via scheduled runner
# Applies updated price/profit calculations to several leve rows across
# the crafting-class worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW).

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# --- ALC ---
Set-LeveRow "ALC" 129 @{
    "H" = 936.4400000000001
    "I" = 350
    "J" = 948.40814
    "K" = 1050
    "L" = 2845.22442
    "M" = 3950
    "N" = -12845.22442
}

Set-LeveRow "ALC" 134 @{
    "H" = 39937.855
    "J" = 39937.855
    "L" = 39937.855
    "N" = -50077.855
}

# --- ARM ---
Set-LeveRow "ARM" 2 @{
    "H" = 1814.909
    "I" = 1963.8948
    "J" = 871.3333
    "K" = 1963.8948
    "L" = 871.3333
    "M" = -1850.8948
    "N" = -1097.3333
}

Set-LeveRow "ARM" 61 @{
    "H" = 2302.0476
    "I" = 1788.4706
    "J" = 4484.75
    "K" = 1788.4706
    "L" = 4484.75
    "M" = -1576.4706
    "N" = -4908.75
}

Set-LeveRow "ARM" 102 @{
    "H" = 48381.6
    "J" = 53335.223
    "L" = 53335.223
    "N" = -56579.223
}

Set-LeveRow "ARM" 110 @{
    "H" = 1501.2439
    "I" = 1431.1177
    "J" = 1841.8572
    "K" = 1431.1177
    "L" = 1841.8572
    "M" = 613.8823
    "N" = -5931.8572
}

Set-LeveRow "ARM" 116 @{
    "H" = 1814.909
    "I" = 1963.8948
    "J" = 871.3333
    "K" = 1963.8948
    "L" = 871.3333
    "M" = 330.1052
    "N" = -5459.3333
}

# Row 122: HQ price/leve price drop to 0, so the HQ profit cell (N) becomes
# blank while the NQ profit (M) is updated instead.
$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H122").Value = 1600
$wsARM.Range("I122").Value = 1600
$wsARM.Range("J122").Value = 0
$wsARM.Range("K122").Value = 4800
$wsARM.Range("L122").Value = 0
$wsARM.Range("N122").Value = ""
$wsARM.Range("M122").Value = -2350

Set-LeveRow "ARM" 132 @{
    "H" = 19233000
    "I" = 27779178
    "J" = 4099.5
    "K" = 83337534
    "L" = 12298.5
    "M" = -83335004
    "N" = -17358.5
}

Set-LeveRow "ARM" 136 @{
    "H" = 2302.0476
    "I" = 1788.4706
    "J" = 4484.75
    "K" = 5365.4118
    "L" = 13454.25
    "M" = -2815.4118
    "N" = -18554.25
}

# --- BSM ---
Set-LeveRow "BSM" 3 @{
    "H" = 1814.909
    "I" = 1963.8948
    "J" = 871.3333
    "K" = 1963.8948
    "L" = 871.3333
    "M" = -1849.8948
    "N" = -1099.3333
}

Set-LeveRow "BSM" 94 @{
    "H" = 1089.409
    "I" = 1050.8096
    "J" = 1900
    "K" = 1050.8096
    "L" = 1900
    "M" = -599.8096
    "N" = -2802
}

Set-LeveRow "BSM" 99 @{
    "H" = 2387.5518
    "I" = 2019.0869
    "J" = 3800
    "K" = 2019.0869
    "L" = 3800
    "M" = -521.0869
    "N" = -6796
}

Set-LeveRow "BSM" 107 @{
    "H" = 1926.3103
    "I" = 1552.5
    "J" = 2757
    "K" = 1552.5
    "L" = 2757
    "M" = 367.5
    "N" = -6597
}

Set-LeveRow "BSM" 134 @{
    "H" = 3306.5557
    "I" = 2703.0476
    "J" = 3608.3096
    "K" = 8109.1428
    "L" = 10824.9288
    "M" = -5574.1428
    "N" = -15894.9288
}

# --- CRP ---
Set-LeveRow "CRP" 99 @{
    "H" = 3766.6667
    "I" = 1500
    "J" = 4900
    "K" = 1500
    "L" = 4900
    "M" = -2
    "N" = -7896
}

Set-LeveRow "CRP" 122 @{
    "H" = 241040
    "I" = 301000
    "J" = 1200
    "K" = 903000
    "L" = 3600
    "M" = -900550
    "N" = -8500
}

Set-LeveRow "CRP" 126 @{
    "H" = 3766.6667
    "I" = 1500
    "J" = 4900
    "K" = 4500
    "L" = 14700
    "M" = -2030
    "N" = -19640
}

# --- CUL ---
Set-LeveRow "CUL" 113 @{
    "H" = 2917.8372
    "I" = 7128.3335
    "J" = 662.2143
    "K" = 21385.0005
    "L" = 1986.6429
    "M" = -19215.0005
    "N" = -6326.6429
}

Set-LeveRow "CUL" 134 @{
    "H" = 31216346
    "I" = 51503412
    "J" = 5474.231
    "K" = 154510236
    "L" = 16422.693
    "M" = -154505166
    "N" = -26562.693
}

# --- GSM ---
Set-LeveRow "GSM" 97 @{
    "H" = 4359.048
    "I" = 2902.6667
    "J" = 8000
    "K" = 2902.6667
    "L" = 8000
    "M" = -2406.6667
    "N" = -8992
}

Set-LeveRow "GSM" 122 @{
    "H" = 1330.375
    "I" = 1028.4286
    "J" = 1565.2222
    "K" = 3085.2858
    "L" = 4695.6666
    "M" = -635.2857999999997
    "N" = -9595.6666
}

Set-LeveRow "GSM" 132 @{
    "H" = 25002702
    "I" = 38463164
    "J" = 4699
    "K" = 115389492
    "L" = 14097
    "M" = -115386962
    "N" = -19157
}

# --- LTW ---
Set-LeveRow "LTW" 61 @{
    "H" = 6502
    "I" = 6502
    "K" = 6502
    "M" = -6300
}

Set-LeveRow "LTW" 100 @{
    "H" = 2349
    "I" = 2284.5715
    "K" = 2284.5715
    "M" = -1743.5715
}

Set-LeveRow "LTW" 113 @{
    "H" = 6502
    "I" = 6502
    "K" = 6502
    "M" = -4332
}
